$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.771.79'
$ws.Range('E2').Value = '  +0.50%  '

# Row 3
$ws.Range('D3').Value = '2.532.92'
$ws.Range('E3').Value = '  +0.85%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.93'
$ws.Range('E5').Value = '  -0.05%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.76'
$ws.Range('E6').Value = '  -1.57%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  -1.16%  '

# Row 9
$ws.Range('D9').Value = '2.532.04'
$ws.Range('E9').Value = '  +0.93%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  -2.36%  '

# Row 11
$ws.Range('E11').Value = '  +1.27%  '

# Row 12
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.344'
$ws.Range('E12').Value = '  +0.28%  '

# Row 13
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.02'
$ws.Range('E13').Value = '  -2.62%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.17'
$ws.Range('E14').Value = '  -2.14%  '

# Row 15
$ws.Range('D15').Value = '2.999.43'
$ws.Range('E15').Value = '  +1.07%  '

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000175'
$ws.Range('E16').Value = '  -1.59%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '67.615.07'
$ws.Range('E17').Value = '  +0.44%  '

# Row 18
$ws.Range('D18').Value = '2.532.20'
$ws.Range('E18').Value = '  +1.77%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.83'
$ws.Range('E19').Value = '  +3.82%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.00'
$ws.Range('E20').Value = '  -0.01%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.43'
$ws.Range('E21').Value = '  +2.19%  '

# Row 22
$ws.Range('E22').Value = '  +62.30%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.16'
$ws.Range('E23').Value = '  -0.70%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.52'
$ws.Range('E24').Value = '  -2.17%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.13'
$ws.Range('E25').Value = '  +1.74%  '

# Row 26
$ws.Range('E26').Value = '  +0.03%  '

# Row 27
$ws.Range('E27').Value = '  -5.13%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.86'
$ws.Range('E28').Value = '  -3.53%  '

# Row 29
$ws.Range('D29').Value = '2.660.82'
$ws.Range('E29').Value = '  +0.15%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0946'
$ws.Range('E30').Value = '  -3.94%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '534.48'
$ws.Range('E31').Value = '  -2.29%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.27'
$ws.Range('E32').Value = '  +0.43%  '

# Row 33
$ws.Range('E33').Value = '  +0.12%  '

# Row 34
$ws.Range('E34').Value = '  -4.46%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.128'
$ws.Range('E35').Value = '  -1.40%  '

# Row 36
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.08%  '

# Row 37
$ws.Range('E37').Value = '  +3.40%  '

# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  -1.70%  '

# Row 39
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.22'
$ws.Range('E39').Value = '  +2.78%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.63'
$ws.Range('E40').Value = '  +0.22%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').Value = '  -1.87%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.06'
$ws.Range('E42').Value = '  -2.06%  '

# Row 43
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.344'
$ws.Range('E43').Value = '  -3.10%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.32%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  -2.79%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.46'
$ws.Range('E46').Value = '  -0.84%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '148.07'
$ws.Range('E47').Value = '  +1.06%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.551'
$ws.Range('E48').Value = '  -1.48%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.69'
$ws.Range('E49').Value = '  -0.63%  '

# Row 50
$ws.Range('D50').Value = '0.0₆0274'
$ws.Range('E50').Value = '  -2.31%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.71'
$ws.Range('E51').Value = '  +1.05%  '
